$wb = $excel.ActiveWorkbook

# --- survey sheet: rename label::language -> label, hint::language -> hint ---
$survey = $wb.Worksheets.Item("survey")
$survey.Range("C1").Value = "label"
$survey.Range("D1").Value = "hint"

# --- choices sheet: rename label::language -> label ---
$choices = $wb.Worksheets.Item("choices")
$choices.Range("C1").Value = "label"

# --- restore/settle active cell selections on survey & choices ---
$choices.Range("A2").Select() | Out-Null
$survey.Range("A2").Select() | Out-Null
